$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (incl. the french date number format) of the last
# existing data row down into the new row, then fill in the new values,
# mirroring how the next day's odometer reading was appended by hand.
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(29, 1).Value = 43763
$ws.Cells.Item(29, 2).Value = 794

# Reflect the new next empty row as the active selection, like Excel does
# after typing a value and pressing Enter.
$ws.Range("B30").Select()
